$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data: the email and result changed
$ws.Range("A2").Value = "qa@impactanalytics.co"
$ws.Range("C2").Value = "pass"

# Update the active selection to B8
$ws.Range("B8").Select()
